# Apply the "scenarios code" edit:
#  - sgen sheet: clear the sample data row (row 2), keep A2's header-ish style
#  - gen sheet: add a new PV generator row (row 2)

$wb = $excel.ActiveWorkbook

# --- sgen sheet: clear row 2 contents (A2:J2) ---
$sgen = $wb.Worksheets.Item("sgen")
$sgen.Range("A2:J2").ClearContents()
$sgen.Activate()
$sgen.Range("A2:J2").Select()

# --- gen sheet: add new row 2 data for a PV generator ---
$gen = $wb.Worksheets.Item("gen")
$gen.Range("A2").Value = 0
$gen.Range("C2").Value = 2
$gen.Range("D2").Value = 6
$gen.Range("E2").Value = 1.03
$gen.Range("G2").Value = -3
$gen.Range("H2").Value = -3
$gen.Range("I2").Value = 1
$gen.Range("J2").Value = $false
$gen.Range("K2").Value = $true
$gen.Range("M2").Value = "pv"
$gen.Range("B2").Value = "PV-gen"
$gen.Range("A2").Select()

$gen.Activate()
